# Apply a permutation of the "variable" data block (columns D, L, M, N, O, P, Q, R, S, T)
# across data rows 2-18 of the active worksheet. Columns A, B, C, E, F, G, H, I, J, K remain
# unchanged per row; only the rest of each row's content is reshuffled to match the values
# that used to sit in a different row before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: target row number -> source row number (where the "new" values originally lived)
$rowMap = @{
    2  = 14
    3  = 8
    4  = 11
    5  = 18
    6  = 4
    7  = 2
    8  = 5
    9  = 6
    10 = 10
    11 = 12
    12 = 13
    13 = 17
    14 = 3
    15 = 9
    16 = 7
    17 = 15
    18 = 16
}

# Snapshot of the original values for the columns that move, keyed by source row number.
# NOTE: use .Value2 for reads (the .Value getter in this runtime misbehaves), .Value for writes.
$snapshot = @{}
foreach ($srcRow in 2..18) {
    $snapshot[$srcRow] = @{
        D = $ws.Cells.Item($srcRow, 4).Value2   # Fecha
        L = $ws.Cells.Item($srcRow, 12).Value2  # Calidad
        M = $ws.Cells.Item($srcRow, 13).Value2  # Volumen
        N = $ws.Cells.Item($srcRow, 14).Value2  # Precio minimo
        O = $ws.Cells.Item($srcRow, 15).Value2  # Precio maximo
        P = $ws.Cells.Item($srcRow, 16).Value2  # Precio promedio ponderado
        Q = $ws.Cells.Item($srcRow, 17).Value2  # Unidad de comercializacion
        R = $ws.Cells.Item($srcRow, 18).Value2  # Origen
        S = $ws.Cells.Item($srcRow, 19).Value2  # Precio $/Kg
        T = $ws.Cells.Item($srcRow, 20).Value2  # Kg / unidad
    }
}

foreach ($targetRow in 2..18) {
    $srcRow = $rowMap[$targetRow]
    $vals = $snapshot[$srcRow]

    $ws.Cells.Item($targetRow, 4).Value = $vals.D
    $ws.Cells.Item($targetRow, 12).Value = $vals.L
    $ws.Cells.Item($targetRow, 13).Value = $vals.M
    $ws.Cells.Item($targetRow, 14).Value = $vals.N
    $ws.Cells.Item($targetRow, 15).Value = $vals.O
    $ws.Cells.Item($targetRow, 16).Value = $vals.P
    $ws.Cells.Item($targetRow, 17).Value = $vals.Q
    $ws.Cells.Item($targetRow, 18).Value = $vals.R
    $ws.Cells.Item($targetRow, 19).Value = $vals.S
    $ws.Cells.Item($targetRow, 20).Value = $vals.T
}
